$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update recomputed values for rows 4-73 (naive forecaster bugfix)
$ws.Range("B4").Value = -0.2000000000000028
$ws.Range("B5").Value = 0.2999999999999829
$ws.Range("B6").Value = 0.2000000000000028
$ws.Range("B7").Value = 0
$ws.Range("B8").Value = -0.7999999999999972
$ws.Range("B9").Value = 0
$ws.Range("B10").Value = 0.2999999999999829
$ws.Range("B11").Value = 0.4000000000000057
$ws.Range("B12").Value = 0.5999999999999943
$ws.Range("B13").Value = 1.200000000000003
$ws.Range("B14").Value = 2
$ws.Range("B15").Value = 0.7999999999999972
$ws.Range("B16").Value = 1.599999999999994
$ws.Range("B17").Value = 1.799999999999997
$ws.Range("B18").Value = 1
$ws.Range("B19").Value = 0.5999999999999943
$ws.Range("B20").Value = -0.4993864180312784
$ws.Range("B21").Value = -0.2000000000000313
$ws.Range("B22").Value = 0.7063371330579002
$ws.Range("B23").Value = 0.4489068848233728
$ws.Range("B24").Value = 1.408791801231501
$ws.Range("B25").Value = 0.5999999999999943
$ws.Range("B26").Value = 1.011699570515816
$ws.Range("B27").Value = 0.9935161553936211
$ws.Range("B28").Value = 1.302043324436823
$ws.Range("B29").Value = 1.299999999999969
$ws.Range("B30").Value = 0.3957849067177932
$ws.Range("B31").Value = 0.4505051707104855
$ws.Range("B32").Value = 0.59458493635276
$ws.Range("B33").Value = 0.4999999999999432
$ws.Range("B34").Value = 0.6016050630459375
$ws.Range("B35").Value = 0.6001107530880319
$ws.Range("B36").Value = 0.5980962996738413
$ws.Range("B37").Value = 0.5999999999999943
$ws.Range("B38").Value = 0.6486835802838442
$ws.Range("B39").Value = 0.6499971057671843
$ws.Range("B40").Value = 0.600949061571356
$ws.Range("B41").Value = 0.5999999999999943
$ws.Range("B42").Value = 0.9401796321600813
$ws.Range("B43").Value = 0.9025026375028489
$ws.Range("B44").Value = 0.9987092793436005
$ws.Range("B45").Value = 0.8989194962581735
$ws.Range("B46").Value = 0.6
$ws.Range("B47").Value = 0.5479150381202658
$ws.Range("B48").Value = 0.5
$ws.Range("B49").Value = 0.7003888945527734
$ws.Range("B50").Value = 0.6487569224423311
$ws.Range("B51").Value = 0.7003348337051136
$ws.Range("B52").Value = -1
$ws.Range("B53").Value = 0.798184409453512
$ws.Range("B54").Value = 0.7972078981263451
$ws.Range("B55").Value = 0.7006632596227007
$ws.Range("B56").Value = 0.7005260309296233
$ws.Range("B57").Value = 1.501866008463963
$ws.Range("B58").Value = 1.500711092634546
$ws.Range("B59").Value = 1.397766293286139
$ws.Range("B60").Value = 0.9013714134183743
$ws.Range("B61").Value = 0.4038529775848758
$ws.Range("B62").Value = 0.9991743200616696
$ws.Range("B63").Value = 0.6350385067502629
$ws.Range("B64").Value = 0.2134433935270721
$ws.Range("B65").Value = 0.3629264220863746
$ws.Range("B66").Value = 0.3091193133532926
$ws.Range("B67").Value = 0.6338683484835599
$ws.Range("B68").Value = 0.2985013435050234
$ws.Range("B69").Value = 0.2049602221247682
$ws.Range("B70").Value = 0.6878071214384534
$ws.Range("B71").Value = 0.5782185388409715
$ws.Range("B72").Value = 0.5490590963348296
$ws.Range("B73").Value = 0.7477535645244302

# Remove trailing rows 74-82 that no longer exist after the bugfix
$ws.Range("A74:B82").EntireRow.Delete()

Write-Host "Applied naive forecaster bugfix"
